$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This change reflects a new handback run: the file "9bdac964-...md" now has
# a failed handback transform (with an error detail message), while
# "cf5e98ae-...md" reverts to "In Translation". Concretely, the data rows
# for these two files swap places (row 6 <-> row 7) across the Overview,
# zh-cn and de-de sheets, row 6 picking up the new failure status/detail.
# ---------------------------------------------------------------------------

function Set-LinkText($ws, $cellRef, $text) {
    $target = $ws.Range($cellRef)
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Row -eq $target.Row -and $h.Range.Column -eq $target.Column) {
            $h.TextToDisplay = $text
        }
    }
}

# ----------------------- Overview sheet -----------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A6").Value = "9bdac964-426c-4a43-bf00-57683e48f178.md"
$wsOverview.Range("B6").Value = "Handback transform failed"
$wsOverview.Range("C6").Value = "Handback transform failed"
$wsOverview.Range("D6").Value = "2016-17-17 14:17:33"

$wsOverview.Range("A7").Value = "cf5e98ae-0b1b-40ae-9764-07869c05bb5f.md"
$wsOverview.Range("B7").Value = "In Translation"
$wsOverview.Range("C7").Value = "In Translation"
$wsOverview.Range("D7").Value = "2016-12-17 14:12:59"

Set-LinkText $wsOverview "A6" "9bdac964-426c-4a43-bf00-57683e48f178.md"
Set-LinkText $wsOverview "A7" "cf5e98ae-0b1b-40ae-9764-07869c05bb5f.md"

# ----------------------- zh-cn sheet ---------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A6").Value = "9bdac964-426c-4a43-bf00-57683e48f178.md"
$wsZhCn.Range("C6").Value = "Handback transform failed"
$wsZhCn.Range("D6").Value = "9bdac964-426c-4a43-bf00-57683e48f178.7860cd1b2b48fdc29524487063fd4c7f7144fa15.zh-cn.xlf"
$wsZhCn.Range("E6").Value = "2016-03-17 14:17:29"
$wsZhCn.Range("K6").Value = "The handback type mt is not match with handoff type ht."

$wsZhCn.Range("A7").Value = "cf5e98ae-0b1b-40ae-9764-07869c05bb5f.md"
$wsZhCn.Range("C7").Value = "In Translation"
$wsZhCn.Range("D7").Value = "cf5e98ae-0b1b-40ae-9764-07869c05bb5f.253ac65e3fabd5d4c3d44d202594a3c77209f05f.zh-cn.xlf"
$wsZhCn.Range("E7").Value = "2016-03-17 14:12:48"

Set-LinkText $wsZhCn "A6" "9bdac964-426c-4a43-bf00-57683e48f178.md"
Set-LinkText $wsZhCn "D6" "9bdac964-426c-4a43-bf00-57683e48f178.7860cd1b2b48fdc29524487063fd4c7f7144fa15.zh-cn.xlf"
Set-LinkText $wsZhCn "A7" "cf5e98ae-0b1b-40ae-9764-07869c05bb5f.md"
Set-LinkText $wsZhCn "D7" "cf5e98ae-0b1b-40ae-9764-07869c05bb5f.253ac65e3fabd5d4c3d44d202594a3c77209f05f.zh-cn.xlf"

# ----------------------- de-de sheet ---------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A6").Value = "9bdac964-426c-4a43-bf00-57683e48f178.md"
$wsDeDe.Range("C6").Value = "Handback transform failed"
$wsDeDe.Range("D6").Value = "9bdac964-426c-4a43-bf00-57683e48f178.7860cd1b2b48fdc29524487063fd4c7f7144fa15.de-de.xlf"
$wsDeDe.Range("E6").Value = "2016-03-17 14:17:33"
$wsDeDe.Range("K6").Value = "The handback type mt is not match with handoff type ht."

$wsDeDe.Range("A7").Value = "cf5e98ae-0b1b-40ae-9764-07869c05bb5f.md"
$wsDeDe.Range("C7").Value = "In Translation"
$wsDeDe.Range("D7").Value = "cf5e98ae-0b1b-40ae-9764-07869c05bb5f.253ac65e3fabd5d4c3d44d202594a3c77209f05f.de-de.xlf"
$wsDeDe.Range("E7").Value = "2016-03-17 14:12:59"

Set-LinkText $wsDeDe "A6" "9bdac964-426c-4a43-bf00-57683e48f178.md"
Set-LinkText $wsDeDe "D6" "9bdac964-426c-4a43-bf00-57683e48f178.7860cd1b2b48fdc29524487063fd4c7f7144fa15.de-de.xlf"
Set-LinkText $wsDeDe "A7" "cf5e98ae-0b1b-40ae-9764-07869c05bb5f.md"
Set-LinkText $wsDeDe "D7" "cf5e98ae-0b1b-40ae-9764-07869c05bb5f.253ac65e3fabd5d4c3d44d202594a3c77209f05f.de-de.xlf"

$wb.Save()
